$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 11) mirroring the pattern of row 10:
# column A gets the next sequential id (number), columns B-E are left blank
# (empty strings so they serialize as inline strings, matching the existing
# pattern used for row 10).
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
